# Add a new "2022-Q3" quarter:
#  1. Insert a new summary row on "总计" for 2022-Q3 (pushes existing rows down).
#  2. Create a new "2022-Q3" sheet (cloned from "2022-Q2" so sheet-level
#     properties/styles match the other quarter sheets) positioned right
#     before "2022-Q2", and fill it with the new holdings data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert row for 2022-Q3 at the top of the data
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Give the new index cell (A2) the same style as the index column below it
# (bold / centred / bordered) instead of the default Excel inherits.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.03

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 2) New "2022-Q3" sheet with the quarter's fund holdings, placed before
#    "2022-Q2" so tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, ...
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'005351"
$q3.Range("C2").Value = "汇添富行业整合主题混合A"
$q3.Range("D2").Value = "'0.67"
$q3.Range("E2").Value = "'83.14"
$q3.Range("F2").Value = "'4.72"
$q3.Range("G2").Value = "'0.0316"
$q3.Range("H2").Value = 8

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'015191"
$q3.Range("C3").Value = "汇添富行业整合主题混合D"
$q3.Range("D3").Value = "'0.00"
$q3.Range("E3").Value = "'83.14"
$q3.Range("F3").Value = "'4.72"
$q3.Range("G3").Value = 0
$q3.Range("H3").Value = 8

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'015190"
$q3.Range("C4").Value = "汇添富行业整合主题混合C"
$q3.Range("D4").Value = "'0.00"
$q3.Range("E4").Value = "'83.14"
$q3.Range("F4").Value = "'4.72"
$q3.Range("G4").Value = 0
$q3.Range("H4").Value = 8

# Restore the originally-active tab ("2021-Q1" was tabSelected before this
# edit); creating/copying sheets above moved the active tab to the new one.
$wb.Worksheets.Item("2021-Q1").Activate()
